$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellValues = @{
    "B2" = 0.2074829931972789
    "C2" = 0.5578231292517006
    "J2" = 0.01020408163265306
    "P2" = 0.1530612244897959
    "S2" = 0.07142857142857142
    "B3" = 0.005882352941176471
    "C3" = 0.02941176470588235
    "J3" = 0.01176470588235294
    "P3" = 0.8117647058823529
    "S3" = 0.1411764705882353
    "J4" = 0.02380952380952381
    "P4" = 0.7142857142857143
    "S4" = 0.2619047619047619
    "S5" = 1
    "B6" = 0.07234042553191489
    "D6" = 0.008510638297872341
    "E6" = 0.008510638297872341
    "F6" = 0.05531914893617021
    "J6" = 0.2723404255319149
    "O6" = 0.01276595744680851
    "Q6" = 0.1702127659574468
    "R6" = 0.0851063829787234
    "S6" = 0.3148936170212766
    "B7" = 0.06201550387596899
    "D7" = 0.01550387596899225
    "F7" = 0.08527131782945736
    "J7" = 0.1705426356589147
    "O7" = 0.02325581395348837
    "Q7" = 0.2093023255813954
    "R7" = 0.04651162790697674
    "S7" = 0.3875968992248062
    "B8" = 0.09429280397022333
    "D8" = 0.02481389578163772
    "E8" = 0.002481389578163772
    "F8" = 0.06699751861042183
    "J8" = 0.1290322580645161
    "O8" = 0.02233250620347394
    "Q8" = 0.1861042183622829
    "R8" = 0.109181141439206
    "S8" = 0.3647642679900744
    "B9" = 0.1171171171171171
    "D9" = 0.01801801801801802
    "F9" = 0.06306306306306306
    "J9" = 0.1216216216216216
    "O9" = 0.02702702702702703
    "Q9" = 0.2027027027027027
    "R9" = 0.1261261261261261
    "S9" = 0.3243243243243243
    "B10" = 0.1217105263157895
    "D10" = 0.01973684210526316
    "E10" = 0.0008223684210526315
    "F10" = 0.07648026315789473
    "J10" = 0.1151315789473684
    "O10" = 0.01973684210526316
    "Q10" = 0.21875
    "R10" = 0.08141447368421052
    "S10" = 0.3462171052631579
    "G11" = 0.1509433962264151
    "J11" = 0.1273584905660377
    "K11" = 0.2028301886792453
    "L11" = 0.5094339622641509
    "S11" = 0.009433962264150943
    "G12" = 0.6779661016949152
    "J12" = 0.211864406779661
    "K12" = 0.008474576271186441
    "L12" = 0.07627118644067797
    "S12" = 0.02542372881355932
    "G13" = 0.6666666666666666
    "J13" = 0.3
    "S13" = 0.03333333333333333
    "G14" = 1
    "F15" = 0.02362204724409449
    "H15" = 0.1692913385826772
    "I15" = 0.1062992125984252
    "J15" = 0.3425196850393701
    "K15" = 0.02755905511811024
    "M15" = 0.01181102362204724
    "N15" = 0.003937007874015748
    "O15" = 0.06299212598425197
    "S15" = 0.2519685039370079
    "F16" = 0.03398058252427184
    "H16" = 0.1601941747572816
    "I16" = 0.0970873786407767
    "J16" = 0.3980582524271845
    "K16" = 0.1019417475728155
    "M16" = 0.03883495145631068
    "O16" = 0.06310679611650485
    "S16" = 0.1067961165048544
    "F17" = 0.02
    "H17" = 0.1955555555555556
    "I17" = 0.1
    "J17" = 0.4222222222222222
    "K17" = 0.06666666666666667
    "M17" = 0.008888888888888889
    "N17" = 0.002222222222222222
    "O17" = 0.08888888888888889
    "S17" = 0.09555555555555556
    "F18" = 0.005050505050505051
    "H18" = 0.1919191919191919
    "I18" = 0.1262626262626263
    "J18" = 0.5050505050505051
    "K18" = 0.07575757575757576
    "N18" = 0.005050505050505051
    "O18" = 0.02525252525252525
    "S18" = 0.06565656565656566
    "F19" = 0.02720450281425891
    "H19" = 0.1857410881801126
    "I19" = 0.09380863039399624
    "J19" = 0.3874296435272045
    "K19" = 0.08536585365853659
    "M19" = 0.01594746716697936
    "N19" = 0.0009380863039399625
    "O19" = 0.0975609756097561
    "S19" = 0.1060037523452158
}

foreach ($cell in $cellValues.Keys) {
    $ws.Range($cell).Value = $cellValues[$cell]
}

Write-Output "Updated $($cellValues.Count) cells"
